$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H28").Value = 1664.5
$ws_ALC.Range("I28").Value = 874.75
$ws_ALC.Range("K28").Value = 874.75
$ws_ALC.Range("M28").Value = -389.75
$ws_ALC.Range("H33").Value = 622.25
$ws_ALC.Range("I33").Value = 741.6
$ws_ALC.Range("K33").Value = 741.6
$ws_ALC.Range("M33").Value = -512.6
$ws_ALC.Range("H53").Value = 3633.9285
$ws_ALC.Range("J53").Value = 4609.2
$ws_ALC.Range("L53").Value = 4609.2
$ws_ALC.Range("N53").Value = -5883.2
$ws_ALC.Range("H64").Value = 18003
$ws_ALC.Range("J64").Value = 18003
$ws_ALC.Range("L64").Value = 18003
$ws_ALC.Range("N64").Value = -18499
$ws_ALC.Range("H67").Value = 18003
$ws_ALC.Range("J67").Value = 18003
$ws_ALC.Range("L67").Value = 18003
$ws_ALC.Range("N67").Value = -19719
$ws_ALC.Range("H111").Value = 1419.6
$ws_ALC.Range("I111").Value = 1449
$ws_ALC.Range("K111").Value = 4347
$ws_ALC.Range("M111").Value = -1280
$ws_ALC.Range("H113").Value = 7909.0835
$ws_ALC.Range("I113").Value = 8831.666999999999
$ws_ALC.Range("J113").Value = 7601.5557
$ws_ALC.Range("K113").Value = 8831.666999999999
$ws_ALC.Range("L113").Value = 7601.5557
$ws_ALC.Range("M113").Value = -5577.666999999999
$ws_ALC.Range("N113").Value = -14109.5557
$ws_ALC.Range("H115").Value = 1345
$ws_ALC.Range("I115").Value = 267.5
$ws_ALC.Range("J115").Value = 3500
$ws_ALC.Range("K115").Value = 802.5
$ws_ALC.Range("L115").Value = 10500
$ws_ALC.Range("M115").Value = 764.5
$ws_ALC.Range("N115").Value = -13634

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 10243.917
$ws_ARM.Range("I2").Value = 825.2222
$ws_ARM.Range("K2").Value = 825.2222
$ws_ARM.Range("M2").Value = -712.2222
$ws_ARM.Range("H32").Value = 4923.246
$ws_ARM.Range("I32").Value = 3689.8794
$ws_ARM.Range("J32").Value = 15142.571
$ws_ARM.Range("K32").Value = 3689.8794
$ws_ARM.Range("L32").Value = 15142.571
$ws_ARM.Range("M32").Value = -3402.8794
$ws_ARM.Range("N32").Value = -15716.571
$ws_ARM.Range("H61").Value = 8853.380999999999
$ws_ARM.Range("I61").Value = 7127.1333
$ws_ARM.Range("K61").Value = 7127.1333
$ws_ARM.Range("M61").Value = -6915.1333
$ws_ARM.Range("H74").Value = 23813624
$ws_ARM.Range("I74").Value = 55558804
$ws_ARM.Range("J74").Value = 4738.25
$ws_ARM.Range("K74").Value = 55558804
$ws_ARM.Range("L74").Value = 4738.25
$ws_ARM.Range("M74").Value = -55557930
$ws_ARM.Range("N74").Value = -6486.25
$ws_ARM.Range("H77").Value = 23813624
$ws_ARM.Range("I77").Value = 55558804
$ws_ARM.Range("J77").Value = 4738.25
$ws_ARM.Range("K77").Value = 277794020
$ws_ARM.Range("L77").Value = 23691.25
$ws_ARM.Range("M77").Value = -277789652
$ws_ARM.Range("N77").Value = -32427.25
$ws_ARM.Range("H116").Value = 10243.917
$ws_ARM.Range("I116").Value = 825.2222
$ws_ARM.Range("K116").Value = 825.2222
$ws_ARM.Range("M116").Value = 1468.7778
$ws_ARM.Range("H122").Value = 3394.5
$ws_ARM.Range("I122").Value = 2851.5833
$ws_ARM.Range("K122").Value = 8554.749899999999
$ws_ARM.Range("M122").Value = -6104.749899999999
$ws_ARM.Range("H132").Value = 2731.9607
$ws_ARM.Range("I132").Value = 2024.7906
$ws_ARM.Range("K132").Value = 6074.3718
$ws_ARM.Range("M132").Value = -3544.3718
$ws_ARM.Range("H136").Value = 8853.380999999999
$ws_ARM.Range("I136").Value = 7127.1333
$ws_ARM.Range("K136").Value = 21381.3999
$ws_ARM.Range("M136").Value = -18831.3999

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 10243.917
$ws_BSM.Range("I3").Value = 825.2222
$ws_BSM.Range("K3").Value = 825.2222
$ws_BSM.Range("M3").Value = -711.2222
$ws_BSM.Range("H134").Value = 2542.625
$ws_BSM.Range("I134").Value = 1547.7142
$ws_BSM.Range("K134").Value = 4643.142599999999
$ws_BSM.Range("M134").Value = -2108.142599999999

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H58").Value = 2577.2593
$ws_CRP.Range("I58").Value = 1298.8636
$ws_CRP.Range("K58").Value = 1298.8636
$ws_CRP.Range("M58").Value = -1095.8636
$ws_CRP.Range("H136").Value = 2577.2593
$ws_CRP.Range("I136").Value = 1298.8636
$ws_CRP.Range("K136").Value = 3896.5908
$ws_CRP.Range("M136").Value = -1346.5908

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H7").Value = 318.18182
$ws_CUL.Range("I7").Value = 381.57144
$ws_CUL.Range("J7").Value = 207.25
$ws_CUL.Range("K7").Value = 1144.71432
$ws_CUL.Range("L7").Value = 621.75
$ws_CUL.Range("M7").Value = -1032.71432
$ws_CUL.Range("N7").Value = -845.75
$ws_CUL.Range("H9").Value = 110473.336
$ws_CUL.Range("I9").Value = 237600.5
$ws_CUL.Range("J9").Value = 8771.6
$ws_CUL.Range("K9").Value = 712801.5
$ws_CUL.Range("L9").Value = 26314.8
$ws_CUL.Range("M9").Value = -712577.5
$ws_CUL.Range("N9").Value = -26762.8
$ws_CUL.Range("H92").Value = 6426
$ws_CUL.Range("I92").Value = 3555
$ws_CUL.Range("J92").Value = 7000.2
$ws_CUL.Range("K92").Value = 10665
$ws_CUL.Range("L92").Value = 21000.6
$ws_CUL.Range("N92").Value = -23496.6
$ws_CUL.Range("H132").Value = 3233.3704
$ws_CUL.Range("J132").Value = 3703.8462
$ws_CUL.Range("L132").Value = 33334.6158
$ws_CUL.Range("N132").Value = -38394.6158

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 1837.9667
$ws_GSM.Range("I102").Value = 975.6087
$ws_GSM.Range("K102").Value = 975.6087
$ws_GSM.Range("M102").Value = 646.3913
$ws_GSM.Range("H122").Value = 7657.9644
$ws_GSM.Range("J122").Value = 9251.923000000001
$ws_GSM.Range("L122").Value = 27755.769
$ws_GSM.Range("N122").Value = -32655.769
$ws_GSM.Range("H126").Value = 5020.231
$ws_GSM.Range("I126").Value = 2499.75
$ws_GSM.Range("J126").Value = 6140.4443
$ws_GSM.Range("K126").Value = 7499.25
$ws_GSM.Range("L126").Value = 18421.3329
$ws_GSM.Range("M126").Value = -5029.25
$ws_GSM.Range("N126").Value = -23361.3329
$ws_GSM.Range("H132").Value = 3279.12
$ws_GSM.Range("I132").Value = 2671.875
$ws_GSM.Range("J132").Value = 4358.6665
$ws_GSM.Range("K132").Value = 8015.625
$ws_GSM.Range("L132").Value = 13075.9995
$ws_GSM.Range("M132").Value = -5485.625
$ws_GSM.Range("N132").Value = -18135.9995
$ws_GSM.Range("H135").Value = 68699.92999999999
$ws_GSM.Range("J135").Value = 68699.92999999999
$ws_GSM.Range("L135").Value = 68699.92999999999
$ws_GSM.Range("N135").Value = -78839.92999999999

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 2579.2424
$ws_LTW.Range("I22").Value = 1504.826
$ws_LTW.Range("K22").Value = 1504.826
$ws_LTW.Range("M22").Value = -1209.826
$ws_LTW.Range("H27").Value = 2579.2424
$ws_LTW.Range("I27").Value = 1504.826
$ws_LTW.Range("K27").Value = 1504.826
$ws_LTW.Range("M27").Value = -1397.826
$ws_LTW.Range("H46").Value = 6597.485
$ws_LTW.Range("I46").Value = 4760.273
$ws_LTW.Range("K46").Value = 4760.273
$ws_LTW.Range("M46").Value = -4572.273
$ws_LTW.Range("H55").Value = 1564451.6
$ws_LTW.Range("I55").Value = 2382041.5
$ws_LTW.Range("K55").Value = 2382041.5
$ws_LTW.Range("M55").Value = -2381868.5
$ws_LTW.Range("H56").Value = 54166.668
$ws_LTW.Range("I56").Value = 62500
$ws_LTW.Range("K56").Value = 62500
$ws_LTW.Range("M56").Value = -61809
$ws_LTW.Range("H61").Value = 7243.8667
$ws_LTW.Range("I61").Value = 6012.8184
$ws_LTW.Range("K61").Value = 6012.8184
$ws_LTW.Range("M61").Value = -5810.8184
$ws_LTW.Range("H113").Value = 7243.8667
$ws_LTW.Range("I113").Value = 6012.8184
$ws_LTW.Range("K113").Value = 6012.8184
$ws_LTW.Range("M113").Value = -3842.8184
$ws_LTW.Range("H132").Value = 5915.825
$ws_LTW.Range("I132").Value = 6749.5557
$ws_LTW.Range("J132").Value = 4184.231
$ws_LTW.Range("K132").Value = 20248.6671
$ws_LTW.Range("L132").Value = 12552.693
$ws_LTW.Range("M132").Value = -17718.6671
$ws_LTW.Range("N132").Value = -17612.693
$ws_LTW.Range("H136").Value = 3482.361
$ws_LTW.Range("I136").Value = 2238
$ws_LTW.Range("K136").Value = 6714
$ws_LTW.Range("M136").Value = -4164
$ws_LTW.Range("H141").Value = 81342.28999999999
$ws_LTW.Range("J141").Value = 81342.28999999999
$ws_LTW.Range("L141").Value = 81342.28999999999
$ws_LTW.Range("N141").Value = -91702.28999999999

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H51").Value = 29997
$ws_WVR.Range("I51").Value = 29997
$ws_WVR.Range("K51").Value = 29997
$ws_WVR.Range("M51").Value = -29487
$ws_WVR.Range("H107").Value = 1666.091
$ws_WVR.Range("I107").Value = 1576.7142
$ws_WVR.Range("K107").Value = 4730.142599999999
$ws_WVR.Range("M107").Value = -2810.142599999999
$ws_WVR.Range("H136").Value = 2758.5095
$ws_WVR.Range("I136").Value = 1867.3846
$ws_WVR.Range("K136").Value = 5602.1538
$ws_WVR.Range("M136").Value = -3052.1538
$ws_WVR.Range("H137").Value = 69794
$ws_WVR.Range("J137").Value = 69794
$ws_WVR.Range("L137").Value = 69794
$ws_WVR.Range("N137").Value = -79994

# New cell addition: CUL!M92
$ws_CUL.Range("M92").Value = -9417
